$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E range to text format so numeric-looking strings are preserved exactly
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range('D2').Value2 = '28.964.57'
$ws.Range('E2').Value2 = '  -0.18%  '

$ws.Range('D3').Value2 = '1.826.05'
$ws.Range('E3').Value2 = '  -0.29%  '

$ws.Range('D4').Value2 = '0.9945'
$ws.Range('E4').Value2 = '  -0.53%  '

$ws.Range('D5').Value2 = '243.73'
$ws.Range('E5').Value2 = '  +0.92%  '

$ws.Range('D6').Value2 = '0.6303'
$ws.Range('E6').Value2 = '  +0.75%  '

$ws.Range('D7').Value2 = '0.9975'
$ws.Range('E7').Value2 = '  -0.33%  '

$ws.Range('D8').Value2 = '0.07489'
$ws.Range('E8').Value2 = '  -1.66%  '

$ws.Range('E9').Value2 = '  +1.21%  '

$ws.Range('D10').Value2 = '23.10'
$ws.Range('E10').Value2 = '  +1.84%  '

$ws.Range('D11').Value2 = '0.07690'
$ws.Range('E11').Value2 = '  -0.77%  '

$ws.Range('D12').Value2 = '1.826.36'
$ws.Range('E12').Value2 = '  -0.23%  '

$ws.Range('D13').Value2 = '4.995'
$ws.Range('E13').Value2 = '  +0.93%  '

$ws.Range('D14').Value2 = '0.6666'
$ws.Range('E14').Value2 = '  +0.75%  '

$ws.Range('D15').Value2 = '83.06'
$ws.Range('E15').Value2 = '  +0.95%  '

$ws.Range('D16').Value2 = '0.000009702'
$ws.Range('E16').Value2 = '  +2.64%  '

$ws.Range('D17').Value2 = '6.042'
$ws.Range('E17').Value2 = '  +1.33%  '

$ws.Range('D18').Value2 = '28.999.20'

$ws.Range('D19').Value2 = '12.57'
$ws.Range('E19').Value2 = '  +2.24%  '

$ws.Range('D20').Value2 = '225.10'

$ws.Range('D21').Value2 = '0.9975'
$ws.Range('E21').Value2 = '  -0.29%  '

$ws.Range('D22').Value2 = '7.122'
$ws.Range('E22').Value2 = '  -0.82%  '

$ws.Range('D23').Value2 = '0.9967'
$ws.Range('E23').Value2 = '  -0.46%  '

$ws.Range('D24').Value2 = '159.56'
$ws.Range('E24').Value2 = '  -0.12%  '

$ws.Range('D25').Value2 = '0.1415'
$ws.Range('E25').Value2 = '  +4.14%  '

$ws.Range('D26').Value2 = '8.504'
$ws.Range('E26').Value2 = '  +1.22%  '

$ws.Range('D27').Value2 = '17.89'
$ws.Range('E27').Value2 = '  +0.50%  '

$ws.Range('D28').Value2 = '1.498'
$ws.Range('E28').Value2 = '  +0.38%  '

$ws.Range('D29').Value2 = '4.126'
$ws.Range('E29').Value2 = '  +1.72%  '

$ws.Range('D30').Value2 = '4.052'
$ws.Range('E30').Value2 = '  +0.79%  '

$ws.Range('D31').Value2 = '0.05465'
$ws.Range('E31').Value2 = '  +5.19%  '

$ws.Range('D32').Value2 = '1.201'
$ws.Range('E32').Value2 = '  +0.25%  '

$ws.Range('D33').Value2 = '1.857'
$ws.Range('E33').Value2 = '  +0.81%  '

$ws.Range('D34').Value2 = '0.7437'
$ws.Range('E34').Value2 = '  +1.60%  '

$ws.Range('D35').Value2 = '1.135'
$ws.Range('E35').Value2 = '  -0.92%  '

$ws.Range('D36').Value2 = '2.611'
$ws.Range('E36').Value2 = '  -3.25%  '

$ws.Range('D37').Value2 = '1.243.36'
$ws.Range('E37').Value2 = '  -1.54%  '

$ws.Range('B38').Value2 = 'FraxShare'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value2 = '6.728'
$ws.Range('E38').Value2 = '  +6.97%  '

$ws.Range('D39').Value2 = '0.01782'
$ws.Range('E39').Value2 = '  -0.07%  '

$ws.Range('B40').Value2 = 'MXToken'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value2 = '2.742'
$ws.Range('E40').Value2 = '  -0.52%  '

$ws.Range('D41').Value2 = '0.8995'
$ws.Range('E41').Value2 = '  +0.94%  '

$ws.Range('D42').Value2 = '0.9982'
$ws.Range('E42').Value2 = '  -0.30%  '

$ws.Range('D43').Value2 = '101.53'
$ws.Range('E43').Value2 = '  +0.10%  '

$ws.Range('D44').Value2 = '1.969.51'
$ws.Range('E44').Value2 = '  -0.42%  '

$ws.Range('B45').Value2 = 'Aave'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value2 = '65.29'
$ws.Range('E45').Value2 = '  +1.45%  '

$ws.Range('B46').Value2 = 'BabyDogeCoin'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value2 = '0.00000000121'
$ws.Range('E46').Value2 = '  +0.76%  '

$ws.Range('D47').Value2 = '0.5062'
$ws.Range('E47').Value2 = '  -0.96%  '

$ws.Range('D48').Value2 = '0.4049'
$ws.Range('E48').Value2 = '  +1.95%  '

$ws.Range('B49').Value2 = 'EnergySwap'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value2 = '9.019'
$ws.Range('E49').Value2 = '  +1.74%  '

$ws.Range('B50').Value2 = 'XinFinNetwork'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').Value2 = '0.07269'
$ws.Range('E50').Value2 = '  +0.67%  '

$ws.Range('D51').Value2 = '1.664'
$ws.Range('E51').Value2 = '  +2.12%  '

# Restore default style (remove the temporary text-format style marker)
$fmtRange.Style = "Normal"